# Facies analysis fix (in progress) - update E column (DOC normalized?) values
# on the "info" sheet, restyle the updated cells to center alignment while
# dropping the old font override, clear a few now-unknown values, and move
# the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

# --- Update values in column E ---------------------------------------
$ws.Range("E3").Value = 1.3
$ws.Range("E4").Value = 1.4
$ws.Range("E5").Value = 1.4
$ws.Range("E6").Value = 2.5
$ws.Range("E7").Value = 15.4

$ws.Range("E10").Value = 3.7
$ws.Range("E11").Value = 3.01
$ws.Range("E12").Value = 8.7
$ws.Range("E13").Value = 7.6
$ws.Range("E14").Value = 6.1

$ws.Range("E16").Value = 9.2

# --- Clear values that are no longer known ----------------------------
$ws.Range("E2").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("E17").ClearContents()

# --- Re-format the updated cells: same border, default font, centered -
# (Excel converges on a new shared cell style here; grab the formatting
# from a cell that already carries it - plain border + default font -
# then layer the horizontal-center alignment on top.)
$ws.Range("B2").Copy()
$ws.Range("E3:E7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E3:E7").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

$ws.Range("B2").Copy()
$ws.Range("E10:E14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E10:E14").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

$ws.Range("B2").Copy()
$ws.Range("E16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E16").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

$excel.CutCopyMode = $false

# --- Move the active selection on the info sheet -----------------------
$ws.Range("I6").Select()
